$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1607142857142857
$ws.Range("C2").Value = 0.6178571428571429
$ws.Range("J2").Value = 0.007142857142857143
$ws.Range("P2").Value = 0.125
$ws.Range("S2").Value = 0.08928571428571429
$ws.Range("C3").Value = 0.005681818181818182
$ws.Range("J3").Value = 0.03409090909090909
$ws.Range("P3").Value = 0.8181818181818182
$ws.Range("S3").Value = 0.1420454545454546
$ws.Range("B6").Value = 0.09090909090909091
$ws.Range("D6").Value = 0.01976284584980237
$ws.Range("F6").Value = 0.09486166007905138
$ws.Range("J6").Value = 0.2806324110671937
$ws.Range("O6").Value = 0.01185770750988142
$ws.Range("Q6").Value = 0.1383399209486166
$ws.Range("R6").Value = 0.03952569169960474
$ws.Range("S6").Value = 0.3241106719367589
$ws.Range("B7").Value = 0.1265822784810127
$ws.Range("D7").Value = 0.03164556962025317
$ws.Range("F7").Value = 0.04430379746835443
$ws.Range("J7").Value = 0.0949367088607595
$ws.Range("O7").Value = 0.01265822784810127
$ws.Range("Q7").Value = 0.2278481012658228
$ws.Range("R7").Value = 0.05696202531645569
$ws.Range("S7").Value = 0.4050632911392405
$ws.Range("B8").Value = 0.106508875739645
$ws.Range("D8").Value = 0.01775147928994083
$ws.Range("F8").Value = 0.07100591715976332
$ws.Range("J8").Value = 0.07889546351084813
$ws.Range("O8").Value = 0.01972386587771203
$ws.Range("Q8").Value = 0.1952662721893491
$ws.Range("R8").Value = 0.08086785009861933
$ws.Range("S8").Value = 0.4299802761341223
$ws.Range("B9").Value = 0.09219858156028368
$ws.Range("D9").Value = 0.02836879432624113
$ws.Range("F9").Value = 0.05673758865248227
$ws.Range("J9").Value = 0.09219858156028368
$ws.Range("O9").Value = 0.02127659574468085
$ws.Range("Q9").Value = 0.2340425531914894
$ws.Range("R9").Value = 0.0851063829787234
$ws.Range("S9").Value = 0.3900709219858156
$ws.Range("B10").Value = 0.1140684410646388
$ws.Range("D10").Value = 0.01996197718631179
$ws.Range("E10").Value = 0.001901140684410646
$ws.Range("F10").Value = 0.06844106463878327
$ws.Range("J10").Value = 0.08460076045627377
$ws.Range("O10").Value = 0.01711026615969582
$ws.Range("Q10").Value = 0.2423954372623574
$ws.Range("R10").Value = 0.07984790874524715
$ws.Range("S10").Value = 0.3716730038022814
$ws.Range("G11").Value = 0.1734317343173432
$ws.Range("J11").Value = 0.07011070110701106
$ws.Range("K11").Value = 0.2361623616236162
$ws.Range("L11").Value = 0.5092250922509225
$ws.Range("S11").Value = 0.01107011070110701
$ws.Range("G12").Value = 0.7142857142857143
$ws.Range("J12").Value = 0.2214285714285714
$ws.Range("L12").Value = 0.02142857142857143
$ws.Range("S12").Value = 0.04285714285714286
$ws.Range("F15").Value = 0.04040404040404041
$ws.Range("H15").Value = 0.1919191919191919
$ws.Range("I15").Value = 0.0707070707070707
$ws.Range("J15").Value = 0.3585858585858586
$ws.Range("K15").Value = 0.0707070707070707
$ws.Range("M15").Value = 0.02525252525252525
$ws.Range("O15").Value = 0.0707070707070707
$ws.Range("S15").Value = 0.1717171717171717
$ws.Range("F16").Value = 0.03349282296650718
$ws.Range("H16").Value = 0.2009569377990431
$ws.Range("I16").Value = 0.09569377990430622
$ws.Range("J16").Value = 0.3971291866028708
$ws.Range("K16").Value = 0.1004784688995215
$ws.Range("M16").Value = 0.01913875598086124
$ws.Range("O16").Value = 0.05263157894736842
$ws.Range("S16").Value = 0.1004784688995215
$ws.Range("F17").Value = 0.03896103896103896
$ws.Range("H17").Value = 0.2294372294372294
$ws.Range("I17").Value = 0.0735930735930736
$ws.Range("J17").Value = 0.4047619047619048
$ws.Range("K17").Value = 0.08874458874458875
$ws.Range("M17").Value = 0.002164502164502165
$ws.Range("O17").Value = 0.0735930735930736
$ws.Range("S17").Value = 0.08874458874458875
$ws.Range("F18").Value = 0.03846153846153846
$ws.Range("H18").Value = 0.2371794871794872
$ws.Range("I18").Value = 0.0641025641025641
$ws.Range("J18").Value = 0.3525641025641026
$ws.Range("K18").Value = 0.08974358974358974
$ws.Range("M18").Value = 0.01923076923076923
$ws.Range("O18").Value = 0.1025641025641026
$ws.Range("S18").Value = 0.09615384615384616
$ws.Range("F19").Value = 0.03584229390681003
$ws.Range("H19").Value = 0.2589605734767025
$ws.Range("I19").Value = 0.05824372759856631
$ws.Range("J19").Value = 0.3261648745519714
$ws.Range("K19").Value = 0.1021505376344086
$ws.Range("M19").Value = 0.02508960573476703
$ws.Range("O19").Value = 0.05824372759856631
$ws.Range("S19").Value = 0.1353046594982079
